# Update scrape timestamp / row-count headers and append the newly
# scraped rows to each of the three sheets (LP1912, LP1912-215, 6203-6173).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "LP1912"  (columns: A Fecha-ish(blank) B Hora_Scrap C Hora_Llegada D Linea E Minutos F Parada G Fecha)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Cells.Item(2,1).Value = "Última actualización: 31/12/2025 16:23:26"
$ws1.Cells.Item(3,1).Value = "Total filas: 1088"

$rows1 = @(
    @("16:23:15", "16:31", "16_P MOR-SANTA ANA", 8,  "LP1912", "31/12/2025"),
    @("16:23:15", "16:33", "23_HERNANDEZ", 10, "LP1912", "31/12/2025"),
    @("16:23:15", "16:36", "17X38_ROMERO", 13, "LP1912", "31/12/2025"),
    @("16:23:15", "16:45", "14_ABASTO", 22, "LP1912", "31/12/2025"),
    @("16:23:15", "16:51", "16_SANTA ANA", 28, "LP1912", "31/12/2025"),
    @("16:23:15", "16:54", "10_OLMOS", 31, "LP1912", "31/12/2025"),
    @("16:23:15", "17:05", "14_ABASTO", 42, "LP1912", "31/12/2025"),
    @("16:23:15", "17:07", "15_ABASTO", 44, "LP1912", "31/12/2025"),
    @("16:23:15", "17:14", "10_OLMOS", 51, "LP1912", "31/12/2025"),
    @("16:23:15", "17:23", "16_SANTA ANA", 60, "LP1912", "31/12/2025"),
    @("16:23:15", "17:25", "11_ETCHEVERRY", 62, "LP1912", "31/12/2025"),
    @("16:23:15", "17:27", "15_ABASTO", 64, "LP1912", "31/12/2025"),
    @("16:23:15", "17:29", "23_HERNANDEZ", 66, "LP1912", "31/12/2025"),
    @("16:23:15", "17:34", "10_OLMOS", 71, "LP1912", "31/12/2025"),
    @("16:23:15", "17:35", "16_P MOR-SANTA ANA", 72, "LP1912", "31/12/2025"),
    @("16:23:15", "17:38", "17X38_ROMERO", 75, "LP1912", "31/12/2025"),
    @("16:23:15", "17:45", "16_SANTA ANA", 82, "LP1912", "31/12/2025"),
    @("16:23:15", "17:51", "215_EL PELIGRO", 88, "LP1912", "31/12/2025")
)

$r = 1072
foreach ($row in $rows1) {
    $ws1.Cells.Item($r,1).Value = ""
    $ws1.Cells.Item($r,2).Value = $row[0]
    $ws1.Cells.Item($r,3).Value = $row[1]
    $ws1.Cells.Item($r,4).Value = $row[2]
    $ws1.Cells.Item($r,5).Value = $row[3]
    $ws1.Cells.Item($r,6).Value = $row[4]
    $ws1.Cells.Item($r,7).Value = $row[5]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Sheet 2: "LP1912-215"  (columns: A(blank) B Fecha C Hora_Scrap D Hora_Llegada E Linea F Minutos G Parada)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2,1).Value = "Última actualización: 31/12/2025 16:23:26"
$ws2.Cells.Item(3,1).Value = "Total filas: 74"

$ws2.Cells.Item(75,1).Value = ""
$ws2.Cells.Item(75,2).Value = "31/12/2025"
$ws2.Cells.Item(75,3).Value = "16:23:15"
$ws2.Cells.Item(75,4).Value = "17:51"
$ws2.Cells.Item(75,5).Value = "215_EL PELIGRO"
$ws2.Cells.Item(75,6).Value = 88
$ws2.Cells.Item(75,7).Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet 3: "6203-6173"  (columns: A(blank) B Fecha C Hora_Scrap D Hora_Llegada E Linea F Minutos G Parada)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(2,1).Value = "Última actualización: 31/12/2025 16:23:26"
$ws3.Cells.Item(3,1).Value = "Total filas: 134"

$ws3.Cells.Item(134,1).Value = ""
$ws3.Cells.Item(134,2).Value = "31/12/2025"
$ws3.Cells.Item(134,3).Value = "16:23:26"
$ws3.Cells.Item(134,4).Value = "16:26"
$ws3.Cells.Item(134,5).Value = "215B_LP-P MOR-1 Y 57"
$ws3.Cells.Item(134,6).Value = 3
$ws3.Cells.Item(134,7).Value = "L6173"

$ws3.Cells.Item(135,1).Value = ""
$ws3.Cells.Item(135,2).Value = "31/12/2025"
$ws3.Cells.Item(135,3).Value = "16:23:21"
$ws3.Cells.Item(135,4).Value = "16:58"
$ws3.Cells.Item(135,5).Value = "215C_LA PLATA"
$ws3.Cells.Item(135,6).Value = 35
$ws3.Cells.Item(135,7).Value = "L6203"
